# The workbook gained a new daily price column ("09-dec") that was inserted
# right before the "01-oct." column block (column EL, the 142nd column) on
# the "Prix Spot" sheet. Inserting a whole column shifts every existing
# column from EL ("01-oct.") through FP ("31-oct.") one position to the
# right (to EM..FQ), matching the diff's dimension change from A1:FP25 to
# A1:FQ25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a brand-new column at EL (column index 142). Everything that used
# to live in EL..FP slides one column right, into EM..FQ.
$ws.Columns.Item(142).Insert()

# Header row: label the freshly inserted column.
$ws.Cells.Item(1, 142).Value = "09-dec"

# Data rows 2-25: the new column has no recorded price for this date, same
# as the other still-missing dates in the sheet, so it gets the "-" marker.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 142).Value = "-"
}
